$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, derived from the author's re-scraped price/volume snapshot.
$edits = @(
    @{Cell="D2"; Value="308.41"},
    @{Cell="E2"; Value="1.22%"},
    @{Cell="D3"; Value="36.35"},
    @{Cell="E3"; Value="1.22%"},
    @{Cell="D4"; Value="5.058"},
    @{Cell="E4"; Value="0.73%"},
    @{Cell="D5"; Value="0.08120"},
    @{Cell="E5"; Value="0.83%"},
    @{Cell="D6"; Value="1.989"},
    @{Cell="E6"; Value="5.80%"},
    @{Cell="E7"; Value="-0.05%"},
    @{Cell="D8"; Value="0.9296"},
    @{Cell="E8"; Value="-0.20%"},
    @{Cell="D9"; Value="0.1471"},
    @{Cell="E9"; Value="12.64%"},
    @{Cell="D10"; Value="0.1933"},
    @{Cell="E10"; Value="1.82%"},
    @{Cell="D11"; Value="0.09099"},
    @{Cell="E11"; Value="-0.92%"},
    @{Cell="D12"; Value="0.03518"},
    @{Cell="E12"; Value="-0.03%"},
    @{Cell="D13"; Value="0.09858"},
    @{Cell="E13"; Value="-0.37%"},
    @{Cell="D14"; Value="0.001420"},
    @{Cell="E14"; Value="0.23%"},
    @{Cell="D15"; Value="0.006310"},
    @{Cell="E15"; Value="-0.30%"},
    @{Cell="D16"; Value="3.851"},
    @{Cell="E16"; Value="6.35%"},
    @{Cell="D17"; Value="4.157"},
    @{Cell="E17"; Value="0.09%"},
    @{Cell="D18"; Value="3.423"},
    @{Cell="E18"; Value="6.43%"},
    @{Cell="D19"; Value="0.3449"},
    @{Cell="E19"; Value="-0.05%"},
    @{Cell="D20"; Value="0.1327"},
    @{Cell="E20"; Value="-0.58%"},
    @{Cell="D21"; Value="4.816"},
    @{Cell="E21"; Value="-7.73%"},
    @{Cell="E22"; Value="-7.45%"},
    @{Cell="D23"; Value="0.04365"},
    @{Cell="E23"; Value="-1.25%"},
    @{Cell="D24"; Value="0.001236"},
    @{Cell="E24"; Value="-0.03%"},
    @{Cell="D25"; Value="0.004159"},
    @{Cell="E25"; Value="-11.64%"},
    @{Cell="D27"; Value="0.0001302"},
    @{Cell="E27"; Value="-0.02%"},
    @{Cell="D39"; Value="0.02063"},
    @{Cell="E39"; Value="6.08%"},
    @{Cell="D40"; Value="0.05114"},
    @{Cell="E40"; Value="-0.89%"},
    @{Cell="D41"; Value="0.007469"},
    @{Cell="E41"; Value="-1.10%"},
    @{Cell="D42"; Value="0.01013"},
    @{Cell="E42"; Value="-0.44%"},
    @{Cell="D43"; Value="0.1369"},
    @{Cell="E43"; Value="-0.06%"},
    @{Cell="D44"; Value="0.002124"},
    @{Cell="E44"; Value="-2.33%"},
    @{Cell="D45"; Value="0.009691"},
    @{Cell="E45"; Value="-10.02%"},
    @{Cell="D46"; Value="0.00006282"},
    @{Cell="E46"; Value="-1.02%"},
    @{Cell="D47"; Value="0.00000000751"},
    @{Cell="E47"; Value="-0.03%"},
    @{Cell="D49"; Value="0.001603"},
    @{Cell="E49"; Value="-3.55%"},
    @{Cell="D50"; Value="0.00002103"},
    @{Cell="E50"; Value="-0.03%"},
    @{Cell="D51"; Value="0.0002003"},
    @{Cell="E51"; Value="-0.03%"}
)

foreach ($edit in $edits) {
    $cell = $ws.Range($edit.Cell)
    # Force text storage so numeric-looking strings ("308.41", "1.22%")
    # round-trip as literal text instead of being parsed into a number/percentage.
    $cell.NumberFormat = "@"
    $cell.Value = $edit.Value
    # Drop back to the workbook default style so no stray number format lingers
    # on the cell (matches the source cells, which carry no explicit style).
    $cell.Style = "Normal"
}
